$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Configuration")

# Insert a new row 2 ("Subscription slug" / "nv1") and push the existing
# "Pulumi Resource Group" / "Pulumi Storage Account" rows down.
$ws1.Rows.Item(2).Insert()
$ws1.Range("A2").Value = "Subscription slug"
$ws1.Range("B2").Value = "nv1"

# Add the new "Deployments" worksheet right after "Configuration".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Deployments"

$ws2.Columns.Item(1).ColumnWidth = 18.7109375
$ws2.Columns.Item(2).ColumnWidth = 16.5703125

$ws2.Range("A1").Value = "Defines"
$ws2.Range("B1").Value = "Resource Group"
$ws2.Range("C1").Value = "Service"
$ws2.Range("D1").Value = "App"
$ws2.Range("E1").Value = "Region"
$ws2.Range("F1").Value = "sku"

$ws2.Range("A2").Value = "Service"
$ws2.Range("B2").Value = "rg-nvdev-uks"
$ws2.Range("C2").Value = "BBB"

$ws2.Range("B2").Select()

# Restore "Configuration" as the active sheet / active selection.
$ws1.Activate()
$ws1.Range("B2").Select()
